# Apply weekly update: insert a new reading (date 44452) at rows 208-209,
# pushing the rest of the historical series down by two rows, and keep the
# final two rows of history (previously 317-318) appended at the end
# (319-320).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank rows right before the current row 210. This shifts
#    the existing rows 210..318 down to 212..320, which is exactly what we
#    need because the whole historical block (from row 210 on) moves down
#    by two positions.
$ws.Rows("210:211").Insert()

# 2. The two freshly inserted rows (210, 211) are blank. Populate them with
#    a copy of what is currently in rows 208 and 209 (the old readings that
#    need to move down by two rows as well).
$ws.Range("A208:R208").Copy($ws.Range("A210:R210"))
$ws.Range("A209:R209").Copy($ws.Range("A211:R211"))

# 3. Now overwrite rows 208 and 209 in place with the new reading
#    (Fecha = 44452) values.
$ws.Range("D208").Value2 = 44452
$ws.Range("J208").Value2 = 1608
$ws.Range("K208").Value2 = 600
$ws.Range("L208").Value2 = 600
$ws.Range("M208").Value2 = 600
$ws.Range("P208").Value2 = 600

$ws.Range("D209").Value2 = 44452
$ws.Range("J209").Value2 = 1500
$ws.Range("K209").Value2 = 500
$ws.Range("L209").Value2 = 500
$ws.Range("M209").Value2 = 500
$ws.Range("P209").Value2 = 500
